# Applies the "structura de decizie si structura switch" edit:
#   - set zoom level from 205% to 145%
#   - change active selection to H13
#   - fill in "săpt. 6" (column H) attendance of 1 for several students

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set week-6 (column H) attendance values for the affected rows.
$ws.Range("H6").Value  = 1
$ws.Range("H9").Value  = 1
$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H21").Value = 1

# Update the view: zoom and the active cell/selection.
$ws.Range("H13").Select()
$excel.ActiveWindow.Zoom = 145
